$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "link_object_02"
$ws.Range("B3").Value = "Image_Annotations"
$ws.Range("E3").Value = "I_001, I_002, I_006, I_037, I_012, I_030, I_032"
$ws.Range("C3").Value = "Images containing annotations"

$ws.Range("D14").Select()
